$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (e.g. "5.17", "0.0992")
# rather than being auto-coerced to numbers by Excel's smart-entry parsing.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.336.22'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.000.89'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.41'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.20'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +5.86%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.524'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.990.43'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.93%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.17'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +8.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.457'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.05%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.77'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.98%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.495.83'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.28'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +7.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.999.53'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '59.298.75'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '431.16'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.68'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.72%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +6.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.13'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.32'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.85'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.47%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.18'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +10.22%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.84'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.77'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.11'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0992'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.95'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.993'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0761'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +12.81%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.14'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.67'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.40%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +6.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '400.20'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0353'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.757.96'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.80%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.253'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '35.77'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +27.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.998'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.27'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.97%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.45'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.83%  '
